# This script reproduces an "automatic update" re-scrape of the YSTAD
# logging-notification sheet: a subset of data rows got re-ordered (the
# underlying web source apparently changed its listing order) and the
# "Förändrad" (last-changed) date in column C advanced by one day for
# every data row.
#
# Rather than trying to re-derive a sort key (the new order does not
# correspond to any simple ascending/descending sort of a single column),
# we snapshot the affected rows' contents and rewrite them in the target
# order, then bump column C everywhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold plain values (not formulas) for the data rows.
$valueCols = @("A","B","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
# Columns that hold HYPERLINK formulas derived from column A's value.
$formulaCols = @("S","T","V","W","X","Y")

function Get-RowSnapshot($row) {
    $snap = @{}
    foreach ($col in $valueCols) {
        $snap[$col] = $ws.Range("$col$row").Value()
    }
    foreach ($col in $formulaCols) {
        $cell = $ws.Range("$col$row")
        $f = $cell.Formula
        if ($f -ne $null -and $f -ne "") {
            $snap[$col] = $f
        } else {
            $snap[$col] = $null
        }
    }
    return $snap
}

function Set-RowFromSnapshot($row, $snap) {
    foreach ($col in $valueCols) {
        $ws.Range("$col$row").Value = $snap[$col]
    }
    foreach ($col in $formulaCols) {
        if ($snap[$col] -ne $null) {
            $ws.Range("$col$row").Formula = $snap[$col]
        } else {
            $ws.Range("$col$row").Value = $null
        }
    }
}

# Rows whose data is permuted among themselves (the destination set of
# row numbers equals the source set, just shuffled).
$rowsToSnapshot = @(5,6,7,8,9,10,14,15,16,17,19,20,23)

$snapshots = @{}
foreach ($r in $rowsToSnapshot) {
    $snapshots[$r] = Get-RowSnapshot $r
}

# new row -> source (old) row
$mapping = @{
    5  = 10
    6  = 5
    7  = 9
    8  = 7
    9  = 6
    10 = 8
    14 = 23
    15 = 19
    16 = 15
    17 = 20
    19 = 16
    20 = 17
    23 = 14
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    Set-RowFromSnapshot $newRow $snapshots[$srcRow]
}

# Bump the "Förändrad" date (column C) by one day for every data row
# (rows 2 through 23).
for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Range("C$row")
    $current = $cell.Value2
    $cell.Value2 = $current + 1
}
